$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.4003700119637017, 0.39940165052053089),
    @(-0.28485883104582754, 0.28199420486197901),
    @(-0.15293517688902369, 0.15207473440422703),
    @(-0.14007473446776331, 0.13930897172374834),
    @(-0.13330897193964653, 0.13178241233950949),
    @(-0.030818268010124328, 0.030806275419313689),
    @(-0.010806275686897848, 0.010798949805430169),
    @(-0.068778123173203198, 0.068381231563635403),
    @(-0.06238123178980004, 0.062043236003151314),
    @(-0.056043236232859783, 0.055993158237612306),
    @(-0.051493158462754707, 0.05141081268392611),
    @(-0.045410812914894905, 0.045157699036439514),
    @(-0.039157699271124891, 0.039088668639872814),
    @(-0.027088668895252077, 0.027055143718514252),
    @(-0.021055143955249989, 0.021028737121962848),
    @(-0.015028737359524591, 0.015004779130276535),
    @(-0.0090047793689151945, 0.008999999751322818),
    @(-0.104542819691126, 0.10442845011085922),
    @(-0.027097319710394974, 0.027014052580276005),
    @(-0.018014052796223368, 0.018004332846574655),
    @(-0.0090043330628253315, 0.0089999997835032985),
    @(-0.12005902819322145, 0.11963775372528218),
    @(-0.11063775394817998, 0.10990848658382291),
    @(-0.042127070625356033, 0.041999999663034693),
    @(-0.094964144551212826, 0.094715461039569959),
    @(-0.088715461262310669, 0.088395782798293965),
    @(-0.082395783022391367, 0.08130535007990769),
    @(-0.075305350309069041, 0.074549678963114374),
    @(-0.062549679215296194, 0.062174558043938788),
    @(-0.042174558323918365, 0.042020904046665919),
    @(-0.027020904312836436, 0.027000972571787685),
    @(-0.0060009728579037969, 0.0059999997625501678)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Target OOXML column width is 14.7109375; ColumnWidth is quantized to 1/6
# character-width increments on save, so 13.833333333333334 is the closest
# input that rounds to the nearest achievable stored width (14.666666666666666).
$ws.Columns.Item(2).ColumnWidth = 13.833333333333334
